$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2049180327868853
$ws.Range("C2").Value = 0.5259562841530054
$ws.Range("J2").Value = 0.01502732240437158
$ws.Range("P2").Value = 0.157103825136612
$ws.Range("S2").Value = 0.09699453551912568
$ws.Range("C3").Value = 0.01518987341772152
$ws.Range("J3").Value = 0.03037974683544304
$ws.Range("P3").Value = 0.7265822784810126
$ws.Range("S3").Value = 0.2278481012658228
$ws.Range("J4").Value = 0.0759493670886076
$ws.Range("P4").Value = 0.6075949367088608
$ws.Range("S4").Value = 0.3164556962025317
$ws.Range("B6").Value = 0.05532786885245902
$ws.Range("D6").Value = 0.02049180327868852
$ws.Range("F6").Value = 0.07377049180327869
$ws.Range("J6").Value = 0.2254098360655738
$ws.Range("O6").Value = 0.03483606557377049
$ws.Range("Q6").Value = 0.1762295081967213
$ws.Range("R6").Value = 0.05532786885245902
$ws.Range("S6").Value = 0.3586065573770492
$ws.Range("B7").Value = 0.1016949152542373
$ws.Range("D7").Value = 0.02179176755447942
$ws.Range("F7").Value = 0.03389830508474576
$ws.Range("J7").Value = 0.1476997578692494
$ws.Range("O7").Value = 0.0314769975786925
$ws.Range("Q7").Value = 0.1912832929782082
$ws.Range("R7").Value = 0.07506053268765134
$ws.Range("S7").Value = 0.3970944309927361
$ws.Range("B8").Value = 0.1241758241758242
$ws.Range("D8").Value = 0.01098901098901099
$ws.Range("E8").Value = 0.001098901098901099
$ws.Range("F8").Value = 0.06373626373626373
$ws.Range("J8").Value = 0.1373626373626374
$ws.Range("O8").Value = 0.01758241758241758
$ws.Range("Q8").Value = 0.1868131868131868
$ws.Range("R8").Value = 0.05934065934065934
$ws.Range("S8").Value = 0.3989010989010989
$ws.Range("B9").Value = 0.09313725490196079
$ws.Range("D9").Value = 0.02696078431372549
$ws.Range("E9").Value = 0.002450980392156863
$ws.Range("F9").Value = 0.07107843137254902
$ws.Range("J9").Value = 0.1397058823529412
$ws.Range("O9").Value = 0.02205882352941177
$ws.Range("Q9").Value = 0.2058823529411765
$ws.Range("R9").Value = 0.08333333333333333
$ws.Range("S9").Value = 0.3553921568627451
$ws.Range("B10").Value = 0.1195164075993091
$ws.Range("D10").Value = 0.0155440414507772
$ws.Range("E10").Value = 0.001036269430051813
$ws.Range("F10").Value = 0.07253886010362694
$ws.Range("J10").Value = 0.1340241796200345
$ws.Range("O10").Value = 0.02936096718480138
$ws.Range("Q10").Value = 0.2200345423143351
$ws.Range("R10").Value = 0.06563039723661486
$ws.Range("S10").Value = 0.3423143350604491
$ws.Range("F11").Value = 0.001404494382022472
$ws.Range("G11").Value = 0.1404494382022472
$ws.Range("J11").Value = 0.101123595505618
$ws.Range("K11").Value = 0.1966292134831461
$ws.Range("L11").Value = 0.5533707865168539
$ws.Range("S11").Value = 0.007022471910112359
$ws.Range("G12").Value = 0.7002398081534772
$ws.Range("J12").Value = 0.2182254196642686
$ws.Range("K12").Value = 0.007194244604316547
$ws.Range("L12").Value = 0.03597122302158273
$ws.Range("S12").Value = 0.03836930455635491
$ws.Range("G13").Value = 0.7068965517241379
$ws.Range("J13").Value = 0.2758620689655172
$ws.Range("S13").Value = 0.01724137931034483
$ws.Range("F15").Value = 0.01639344262295082
$ws.Range("H15").Value = 0.1347905282331512
$ws.Range("I15").Value = 0.08196721311475409
$ws.Range("J15").Value = 0.3715846994535519
$ws.Range("K15").Value = 0.0692167577413479
$ws.Range("M15").Value = 0.00546448087431694
$ws.Range("O15").Value = 0.05282331511839709
$ws.Range("S15").Value = 0.2677595628415301
$ws.Range("F16").Value = 0.01839080459770115
$ws.Range("H16").Value = 0.1747126436781609
$ws.Range("I16").Value = 0.0896551724137931
$ws.Range("J16").Value = 0.4344827586206896
$ws.Range("K16").Value = 0.1241379310344828
$ws.Range("M16").Value = 0.009195402298850575
$ws.Range("O16").Value = 0.04827586206896552
$ws.Range("S16").Value = 0.1011494252873563
$ws.Range("F17").Value = 0.01435406698564593
$ws.Range("H17").Value = 0.1779904306220096
$ws.Range("I17").Value = 0.07751196172248803
$ws.Range("J17").Value = 0.4421052631578947
$ws.Range("K17").Value = 0.09473684210526316
$ws.Range("M17").Value = 0.01244019138755981
$ws.Range("N17").Value = 0.001913875598086124
$ws.Range("O17").Value = 0.08325358851674641
$ws.Range("S17").Value = 0.09569377990430622
$ws.Range("F18").Value = 0.02114803625377644
$ws.Range("H18").Value = 0.1722054380664653
$ws.Range("I18").Value = 0.09063444108761329
$ws.Range("J18").Value = 0.4350453172205438
$ws.Range("K18").Value = 0.1329305135951662
$ws.Range("M18").Value = 0.003021148036253776
$ws.Range("O18").Value = 0.05740181268882175
$ws.Range("S18").Value = 0.08761329305135952
$ws.Range("F19").Value = 0.01713632901751714
$ws.Range("H19").Value = 0.1999238385376999
$ws.Range("I19").Value = 0.08149276466108149
$ws.Range("J19").Value = 0.3747143945163747
$ws.Range("K19").Value = 0.1214775323686215
$ws.Range("M19").Value = 0.01408987052551409
$ws.Range("N19").Value = 0.0003808073115003808
$ws.Range("O19").Value = 0.0753998476770754
$ws.Range("S19").Value = 0.1153846153846154
